# Auto-generated edit script: update cryptos list values per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.395.39"
$ws.Range("E2").Value = "  -4.93%  "
$ws.Range("D3").Value = "2.893.23"
$ws.Range("E3").Value = "  -5.52%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "486.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.79%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.415"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.09%  "
$ws.Range("E9").Value = "  -5.47%  "
$ws.Range("E10").Value = "  -8.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.345"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.84%  "
$ws.Range("D12").Value = "3.376.91"
$ws.Range("E12").Value = "  -5.38%  "
$ws.Range("E13").Value = "  -4.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.71%  "
$ws.Range("E15").Value = "  -9.06%  "
$ws.Range("D16").Value = "55.384.34"
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.68%  "
$ws.Range("D18").Value = "2.887.50"
$ws.Range("E18").Value = "  -6.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.01%  "
$ws.Range("E20").Value = "  -7.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.22%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("E24").Value = "  -5.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "61.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.993"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.158"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.78%  "
$ws.Range("D28").Value = "0.0₃0832"
$ws.Range("E28").Value = "  -13.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.19%  "
$ws.Range("E31").Value = "  -6.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.82%  "
$ws.Range("E33").Value = "  -9.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.56"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.74%  "
$ws.Range("E35").Value = "  -9.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "24.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.43%  "
$ws.Range("E38").Value = "  -9.86%  "
$ws.Range("E39").Value = "  -6.92%  "
$ws.Range("D40").Value = "2.923.32"
$ws.Range("E40").Value = "  -5.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.08%  "
$ws.Range("E44").Value = "  -6.56%  "
$ws.Range("D45").Value = "2.082.76"
$ws.Range("E45").Value = "  -10.83%  "
$ws.Range("E46").Value = "  -10.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.905"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0227"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0832"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.84%  "
